$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Title: "Subtype Analysis" -> "BCBM HTMA438: Subtype Analysis" ---
$titleShape = $s.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = "BCBM HTMA438: Subtype Analysis"

# --- Subtitle: "GeoMx" + " DSP Project" -> "GeoMx" + " Digital " + "Spatial Profiling" ---
$subShape = $s.Shapes.Item(2)
$tr = $subShape.TextFrame.TextRange

# Replace the " DSP Project" run's text with " Digital " (keeps its own run/formatting)
$tailRun = $tr.Characters(6, 13)
$tailRun.Text = " Digital "

# Append a new run "Spatial Profiling" right after it
$tr.InsertAfter("Spatial Profiling") | Out-Null
